$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (rows 2-11) from 3 to 2
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = 2
}

# Update selection to F28 as reflected in the saved workbook
$ws.Range("F28").Select()
